$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder countries: Bielorrusia now ranks above Ucrania and Panama ---
# Row 44 previously was Ucrania, row 45 was Panama, row 46 was Bielorrusia.
# New ranking: row 44 = Bielorrusia (fresh data), row 45 = Ucrania (old row44 data),
# row 46 = Panama (old row45 data).
$ws.Range("A44").Value = "Bielorrusia"
$ws.Range("B44").Value = 4204
$ws.Range("C44").Value = 476
$ws.Range("D44").Value = 203
$ws.Range("E44").Value = 3961
$ws.Range("F44").Value = 65
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 40

$ws.Range("A45").Value = "Ucrania"
$ws.Range("B45").Value = 4161
$ws.Range("C45").Value = 397
$ws.Range("D45").Value = 186
$ws.Range("E45").Value = 3859
$ws.Range("F45").Value = 45
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 116

$ws.Range("A46").Value = "Panama"
$ws.Range("B46").Value = 3751
$ws.Range("D46").Value = 75
$ws.Range("E46").Value = 3573
$ws.Range("F46").Value = 106
$ws.Range("H46").Value = 103

# --- Reorder countries: Afganistan now ranks above Cuba ---
# Row 82 previously was Cuba, row 83 was Afganistan.
# New ranking: row 82 = Afganistan (fresh data), row 83 = Cuba (old row82 data).
$ws.Range("A82").Value = "Afganistan"
$ws.Range("B82").Value = 840
$ws.Range("C82").Value = 56
$ws.Range("D82").Value = 54
$ws.Range("E82").Value = 756
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 5
$ws.Range("H82").Value = 30

$ws.Range("A83").Value = "Cuba"
$ws.Range("B83").Value = 814
$ws.Range("D83").Value = 151
$ws.Range("E83").Value = 639
$ws.Range("F83").Value = 15
$ws.Range("H83").Value = 24

# --- Update statistics (Spain, provincias & other countries) ---
# Row 5: Espana
$ws.Range("B5").Value = 182816
$ws.Range("C5").Value = 2157
$ws.Range("D5").Value = 74797
$ws.Range("E5").Value = 88889
$ws.Range("G5").Value = 318
$ws.Range("H5").Value = 19130

# Row 18: Suiza
$ws.Range("B18").Value = 26359
$ws.Range("C18").Value = 23
$ws.Range("E18").Value = 9705
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 1254

# Row 31: Rumania
$ws.Range("B31").Value = 7707
$ws.Range("C31").Value = 491
$ws.Range("D31").Value = 1357
$ws.Range("E31").Value = 5963
$ws.Range("F31").Value = 243

# Row 90: Libano
$ws.Range("B90").Value = 663
$ws.Range("C90").Value = 5
$ws.Range("E90").Value = 557

# Row 115: Vietnam
$ws.Range("D115").Value = 175
$ws.Range("E115").Value = 93
